# #5: cash & deposit done
# Rework the 存款 (deposit) sheet so it carries the same trailing metadata
# columns (property_category/category/date/legislator_name/legislator_id/
# source_file/index) as the other property sheets, and relabel the first
# five headers from literal bank/type names to schema field names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- header row (row 1): rename the generic field headers ---
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# Copy the header cell formatting (bold font + border, style index 1) onto
# the newly-used header cells G1:M1.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- data rows (2-17): fill the new trailing columns ---
# Copy the plain data-row formatting (style index 2, i.e. the un-decorated
# look used by every other data cell) from an existing data cell onto the
# new G2:M17 block first, before any values/number formats are set below.
$ws.Range("B2").Copy()
$ws.Range("G2:M17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G2:G17").Value = "deposit"
$ws.Range("H2:H17").Value = "normal"
# Force text storage for the date column so it keeps the literal
# "2012-04-26" string instead of being reinterpreted as a date serial.
$ws.Range("I2:I17").NumberFormat = "@"
$ws.Range("I2:I17").Value = "2012-04-26"
$ws.Range("J2:J17").Value = "林鴻池"
$ws.Range("K2:K17").Value = 1340
$ws.Range("L2:L17").Value = "tmpdb4b1"

# M (index) mirrors column A's row index for each row.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($r, 1).Value2
}

# F7 was stored as the text "84440" (because the amount is large); it should
# be the plain number 84440, matching every other amount in column F.
$ws.Cells.Item(7, 6).Value = 84440
